$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update existing forecast values for 17.03.2025 (last 4 quarters revised) ---
$ws.Cells.Item(382, 2).Value = 5900
$ws.Cells.Item(383, 2).Value = 5840
$ws.Cells.Item(384, 2).Value = 5770
$ws.Cells.Item(385, 2).Value = 5710

# --- Step 2: Append a full new day (18.03.2025) of 96 quarter-hour rows (386-481) ---
$ws.Cells.Item(386, 1).Value = 45734
$ws.Cells.Item(386, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(386, 2).Value = 5660
$ws.Cells.Item(386, 3).Value = 1
$ws.Cells.Item(386, 4).Value = "18.03.20251"
$ws.Cells.Item(387, 1).Value = 45734.01041666666
$ws.Cells.Item(387, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(387, 2).Value = 5620
$ws.Cells.Item(387, 3).Value = 2
$ws.Cells.Item(387, 4).Value = "18.03.20252"
$ws.Cells.Item(388, 1).Value = 45734.02083333334
$ws.Cells.Item(388, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(388, 2).Value = 5580
$ws.Cells.Item(388, 3).Value = 3
$ws.Cells.Item(388, 4).Value = "18.03.20253"
$ws.Cells.Item(389, 1).Value = 45734.03125
$ws.Cells.Item(389, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(389, 2).Value = 5540
$ws.Cells.Item(389, 3).Value = 4
$ws.Cells.Item(389, 4).Value = "18.03.20254"
$ws.Cells.Item(390, 1).Value = 45734.04166666666
$ws.Cells.Item(390, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(390, 2).Value = 5520
$ws.Cells.Item(390, 3).Value = 5
$ws.Cells.Item(390, 4).Value = "18.03.20255"
$ws.Cells.Item(391, 1).Value = 45734.05208333334
$ws.Cells.Item(391, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(391, 2).Value = 5510
$ws.Cells.Item(391, 3).Value = 6
$ws.Cells.Item(391, 4).Value = "18.03.20256"
$ws.Cells.Item(392, 1).Value = 45734.0625
$ws.Cells.Item(392, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(392, 2).Value = 5500
$ws.Cells.Item(392, 3).Value = 7
$ws.Cells.Item(392, 4).Value = "18.03.20257"
$ws.Cells.Item(393, 1).Value = 45734.07291666666
$ws.Cells.Item(393, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(393, 2).Value = 5490
$ws.Cells.Item(393, 3).Value = 8
$ws.Cells.Item(393, 4).Value = "18.03.20258"
$ws.Cells.Item(394, 1).Value = 45734.08333333334
$ws.Cells.Item(394, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(394, 2).Value = 5500
$ws.Cells.Item(394, 3).Value = 9
$ws.Cells.Item(394, 4).Value = "18.03.20259"
$ws.Cells.Item(395, 1).Value = 45734.09375
$ws.Cells.Item(395, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(395, 2).Value = 5510
$ws.Cells.Item(395, 3).Value = 10
$ws.Cells.Item(395, 4).Value = "18.03.202510"
$ws.Cells.Item(396, 1).Value = 45734.10416666666
$ws.Cells.Item(396, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(396, 2).Value = 5520
$ws.Cells.Item(396, 3).Value = 11
$ws.Cells.Item(396, 4).Value = "18.03.202511"
$ws.Cells.Item(397, 1).Value = 45734.11458333334
$ws.Cells.Item(397, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(397, 2).Value = 5540
$ws.Cells.Item(397, 3).Value = 12
$ws.Cells.Item(397, 4).Value = "18.03.202512"
$ws.Cells.Item(398, 1).Value = 45734.125
$ws.Cells.Item(398, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(398, 2).Value = 5560
$ws.Cells.Item(398, 3).Value = 13
$ws.Cells.Item(398, 4).Value = "18.03.202513"
$ws.Cells.Item(399, 1).Value = 45734.13541666666
$ws.Cells.Item(399, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(399, 2).Value = 5590
$ws.Cells.Item(399, 3).Value = 14
$ws.Cells.Item(399, 4).Value = "18.03.202514"
$ws.Cells.Item(400, 1).Value = 45734.14583333334
$ws.Cells.Item(400, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(400, 2).Value = 5620
$ws.Cells.Item(400, 3).Value = 15
$ws.Cells.Item(400, 4).Value = "18.03.202515"
$ws.Cells.Item(401, 1).Value = 45734.15625
$ws.Cells.Item(401, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(401, 2).Value = 5660
$ws.Cells.Item(401, 3).Value = 16
$ws.Cells.Item(401, 4).Value = "18.03.202516"
$ws.Cells.Item(402, 1).Value = 45734.16666666666
$ws.Cells.Item(402, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(402, 2).Value = 5720
$ws.Cells.Item(402, 3).Value = 17
$ws.Cells.Item(402, 4).Value = "18.03.202517"
$ws.Cells.Item(403, 1).Value = 45734.17708333334
$ws.Cells.Item(403, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(403, 2).Value = 5790
$ws.Cells.Item(403, 3).Value = 18
$ws.Cells.Item(403, 4).Value = "18.03.202518"
$ws.Cells.Item(404, 1).Value = 45734.1875
$ws.Cells.Item(404, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(404, 2).Value = 5880
$ws.Cells.Item(404, 3).Value = 19
$ws.Cells.Item(404, 4).Value = "18.03.202519"
$ws.Cells.Item(405, 1).Value = 45734.19791666666
$ws.Cells.Item(405, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(405, 2).Value = 6000
$ws.Cells.Item(405, 3).Value = 20
$ws.Cells.Item(405, 4).Value = "18.03.202520"
$ws.Cells.Item(406, 1).Value = 45734.20833333334
$ws.Cells.Item(406, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(406, 2).Value = 6140
$ws.Cells.Item(406, 3).Value = 21
$ws.Cells.Item(406, 4).Value = "18.03.202521"
$ws.Cells.Item(407, 1).Value = 45734.21875
$ws.Cells.Item(407, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(407, 2).Value = 6290
$ws.Cells.Item(407, 3).Value = 22
$ws.Cells.Item(407, 4).Value = "18.03.202522"
$ws.Cells.Item(408, 1).Value = 45734.22916666666
$ws.Cells.Item(408, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(408, 2).Value = 6460
$ws.Cells.Item(408, 3).Value = 23
$ws.Cells.Item(408, 4).Value = "18.03.202523"
$ws.Cells.Item(409, 1).Value = 45734.23958333334
$ws.Cells.Item(409, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(409, 2).Value = 6640
$ws.Cells.Item(409, 3).Value = 24
$ws.Cells.Item(409, 4).Value = "18.03.202524"
$ws.Cells.Item(410, 1).Value = 45734.25
$ws.Cells.Item(410, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(410, 2).Value = 6820
$ws.Cells.Item(410, 3).Value = 25
$ws.Cells.Item(410, 4).Value = "18.03.202525"
$ws.Cells.Item(411, 1).Value = 45734.26041666666
$ws.Cells.Item(411, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(411, 2).Value = 7000
$ws.Cells.Item(411, 3).Value = 26
$ws.Cells.Item(411, 4).Value = "18.03.202526"
$ws.Cells.Item(412, 1).Value = 45734.27083333334
$ws.Cells.Item(412, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(412, 2).Value = 7160
$ws.Cells.Item(412, 3).Value = 27
$ws.Cells.Item(412, 4).Value = "18.03.202527"
$ws.Cells.Item(413, 1).Value = 45734.28125
$ws.Cells.Item(413, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(413, 2).Value = 7300
$ws.Cells.Item(413, 3).Value = 28
$ws.Cells.Item(413, 4).Value = "18.03.202528"
$ws.Cells.Item(414, 1).Value = 45734.29166666666
$ws.Cells.Item(414, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(414, 2).Value = 7410
$ws.Cells.Item(414, 3).Value = 29
$ws.Cells.Item(414, 4).Value = "18.03.202529"
$ws.Cells.Item(415, 1).Value = 45734.30208333334
$ws.Cells.Item(415, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(415, 2).Value = 7490
$ws.Cells.Item(415, 3).Value = 30
$ws.Cells.Item(415, 4).Value = "18.03.202530"
$ws.Cells.Item(416, 1).Value = 45734.3125
$ws.Cells.Item(416, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(416, 2).Value = 7540
$ws.Cells.Item(416, 3).Value = 31
$ws.Cells.Item(416, 4).Value = "18.03.202531"
$ws.Cells.Item(417, 1).Value = 45734.32291666666
$ws.Cells.Item(417, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(417, 2).Value = 7560
$ws.Cells.Item(417, 3).Value = 32
$ws.Cells.Item(417, 4).Value = "18.03.202532"
$ws.Cells.Item(418, 1).Value = 45734.33333333334
$ws.Cells.Item(418, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(418, 2).Value = 7540
$ws.Cells.Item(418, 3).Value = 33
$ws.Cells.Item(418, 4).Value = "18.03.202533"
$ws.Cells.Item(419, 1).Value = 45734.34375
$ws.Cells.Item(419, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(419, 2).Value = 7500
$ws.Cells.Item(419, 3).Value = 34
$ws.Cells.Item(419, 4).Value = "18.03.202534"
$ws.Cells.Item(420, 1).Value = 45734.35416666666
$ws.Cells.Item(420, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(420, 2).Value = 7430
$ws.Cells.Item(420, 3).Value = 35
$ws.Cells.Item(420, 4).Value = "18.03.202535"
$ws.Cells.Item(421, 1).Value = 45734.36458333334
$ws.Cells.Item(421, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(421, 2).Value = 7340
$ws.Cells.Item(421, 3).Value = 36
$ws.Cells.Item(421, 4).Value = "18.03.202536"
$ws.Cells.Item(422, 1).Value = 45734.375
$ws.Cells.Item(422, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(422, 2).Value = 7240
$ws.Cells.Item(422, 3).Value = 37
$ws.Cells.Item(422, 4).Value = "18.03.202537"
$ws.Cells.Item(423, 1).Value = 45734.38541666666
$ws.Cells.Item(423, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(423, 2).Value = 7130
$ws.Cells.Item(423, 3).Value = 38
$ws.Cells.Item(423, 4).Value = "18.03.202538"
$ws.Cells.Item(424, 1).Value = 45734.39583333334
$ws.Cells.Item(424, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(424, 2).Value = 7020
$ws.Cells.Item(424, 3).Value = 39
$ws.Cells.Item(424, 4).Value = "18.03.202539"
$ws.Cells.Item(425, 1).Value = 45734.40625
$ws.Cells.Item(425, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(425, 2).Value = 6920
$ws.Cells.Item(425, 3).Value = 40
$ws.Cells.Item(425, 4).Value = "18.03.202540"
$ws.Cells.Item(426, 1).Value = 45734.41666666666
$ws.Cells.Item(426, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(426, 2).Value = 6820
$ws.Cells.Item(426, 3).Value = 41
$ws.Cells.Item(426, 4).Value = "18.03.202541"
$ws.Cells.Item(427, 1).Value = 45734.42708333334
$ws.Cells.Item(427, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(427, 2).Value = 6730
$ws.Cells.Item(427, 3).Value = 42
$ws.Cells.Item(427, 4).Value = "18.03.202542"
$ws.Cells.Item(428, 1).Value = 45734.4375
$ws.Cells.Item(428, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(428, 2).Value = 6640
$ws.Cells.Item(428, 3).Value = 43
$ws.Cells.Item(428, 4).Value = "18.03.202543"
$ws.Cells.Item(429, 1).Value = 45734.44791666666
$ws.Cells.Item(429, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(429, 2).Value = 6570
$ws.Cells.Item(429, 3).Value = 44
$ws.Cells.Item(429, 4).Value = "18.03.202544"
$ws.Cells.Item(430, 1).Value = 45734.45833333334
$ws.Cells.Item(430, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(430, 2).Value = 6500
$ws.Cells.Item(430, 3).Value = 45
$ws.Cells.Item(430, 4).Value = "18.03.202545"
$ws.Cells.Item(431, 1).Value = 45734.46875
$ws.Cells.Item(431, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(431, 2).Value = 6440
$ws.Cells.Item(431, 3).Value = 46
$ws.Cells.Item(431, 4).Value = "18.03.202546"
$ws.Cells.Item(432, 1).Value = 45734.47916666666
$ws.Cells.Item(432, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(432, 2).Value = 6390
$ws.Cells.Item(432, 3).Value = 47
$ws.Cells.Item(432, 4).Value = "18.03.202547"
$ws.Cells.Item(433, 1).Value = 45734.48958333334
$ws.Cells.Item(433, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(433, 2).Value = 6360
$ws.Cells.Item(433, 3).Value = 48
$ws.Cells.Item(433, 4).Value = "18.03.202548"
$ws.Cells.Item(434, 1).Value = 45734.5
$ws.Cells.Item(434, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(434, 2).Value = 6330
$ws.Cells.Item(434, 3).Value = 49
$ws.Cells.Item(434, 4).Value = "18.03.202549"
$ws.Cells.Item(435, 1).Value = 45734.51041666666
$ws.Cells.Item(435, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(435, 2).Value = 6310
$ws.Cells.Item(435, 3).Value = 50
$ws.Cells.Item(435, 4).Value = "18.03.202550"
$ws.Cells.Item(436, 1).Value = 45734.52083333334
$ws.Cells.Item(436, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(436, 2).Value = 6300
$ws.Cells.Item(436, 3).Value = 51
$ws.Cells.Item(436, 4).Value = "18.03.202551"
$ws.Cells.Item(437, 1).Value = 45734.53125
$ws.Cells.Item(437, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(437, 2).Value = 6310
$ws.Cells.Item(437, 3).Value = 52
$ws.Cells.Item(437, 4).Value = "18.03.202552"
$ws.Cells.Item(438, 1).Value = 45734.54166666666
$ws.Cells.Item(438, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(438, 2).Value = 6320
$ws.Cells.Item(438, 3).Value = 53
$ws.Cells.Item(438, 4).Value = "18.03.202553"
$ws.Cells.Item(439, 1).Value = 45734.55208333334
$ws.Cells.Item(439, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(439, 2).Value = 6340
$ws.Cells.Item(439, 3).Value = 54
$ws.Cells.Item(439, 4).Value = "18.03.202554"
$ws.Cells.Item(440, 1).Value = 45734.5625
$ws.Cells.Item(440, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(440, 2).Value = 6370
$ws.Cells.Item(440, 3).Value = 55
$ws.Cells.Item(440, 4).Value = "18.03.202555"
$ws.Cells.Item(441, 1).Value = 45734.57291666666
$ws.Cells.Item(441, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(441, 2).Value = 6410
$ws.Cells.Item(441, 3).Value = 56
$ws.Cells.Item(441, 4).Value = "18.03.202556"
$ws.Cells.Item(442, 1).Value = 45734.58333333334
$ws.Cells.Item(442, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(442, 2).Value = 6450
$ws.Cells.Item(442, 3).Value = 57
$ws.Cells.Item(442, 4).Value = "18.03.202557"
$ws.Cells.Item(443, 1).Value = 45734.59375
$ws.Cells.Item(443, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(443, 2).Value = 6500
$ws.Cells.Item(443, 3).Value = 58
$ws.Cells.Item(443, 4).Value = "18.03.202558"
$ws.Cells.Item(444, 1).Value = 45734.60416666666
$ws.Cells.Item(444, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(444, 2).Value = 6560
$ws.Cells.Item(444, 3).Value = 59
$ws.Cells.Item(444, 4).Value = "18.03.202559"
$ws.Cells.Item(445, 1).Value = 45734.61458333334
$ws.Cells.Item(445, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(445, 2).Value = 6630
$ws.Cells.Item(445, 3).Value = 60
$ws.Cells.Item(445, 4).Value = "18.03.202560"
$ws.Cells.Item(446, 1).Value = 45734.625
$ws.Cells.Item(446, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(446, 2).Value = 6710
$ws.Cells.Item(446, 3).Value = 61
$ws.Cells.Item(446, 4).Value = "18.03.202561"
$ws.Cells.Item(447, 1).Value = 45734.63541666666
$ws.Cells.Item(447, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(447, 2).Value = 6810
$ws.Cells.Item(447, 3).Value = 62
$ws.Cells.Item(447, 4).Value = "18.03.202562"
$ws.Cells.Item(448, 1).Value = 45734.64583333334
$ws.Cells.Item(448, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(448, 2).Value = 6910
$ws.Cells.Item(448, 3).Value = 63
$ws.Cells.Item(448, 4).Value = "18.03.202563"
$ws.Cells.Item(449, 1).Value = 45734.65625
$ws.Cells.Item(449, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(449, 2).Value = 7020
$ws.Cells.Item(449, 3).Value = 64
$ws.Cells.Item(449, 4).Value = "18.03.202564"
$ws.Cells.Item(450, 1).Value = 45734.66666666666
$ws.Cells.Item(450, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(450, 2).Value = 7140
$ws.Cells.Item(450, 3).Value = 65
$ws.Cells.Item(450, 4).Value = "18.03.202565"
$ws.Cells.Item(451, 1).Value = 45734.67708333334
$ws.Cells.Item(451, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(451, 2).Value = 7250
$ws.Cells.Item(451, 3).Value = 66
$ws.Cells.Item(451, 4).Value = "18.03.202566"
$ws.Cells.Item(452, 1).Value = 45734.6875
$ws.Cells.Item(452, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(452, 2).Value = 7360
$ws.Cells.Item(452, 3).Value = 67
$ws.Cells.Item(452, 4).Value = "18.03.202567"
$ws.Cells.Item(453, 1).Value = 45734.69791666666
$ws.Cells.Item(453, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(453, 2).Value = 7470
$ws.Cells.Item(453, 3).Value = 68
$ws.Cells.Item(453, 4).Value = "18.03.202568"
$ws.Cells.Item(454, 1).Value = 45734.70833333334
$ws.Cells.Item(454, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(454, 2).Value = 7570
$ws.Cells.Item(454, 3).Value = 69
$ws.Cells.Item(454, 4).Value = "18.03.202569"
$ws.Cells.Item(455, 1).Value = 45734.71875
$ws.Cells.Item(455, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(455, 2).Value = 7680
$ws.Cells.Item(455, 3).Value = 70
$ws.Cells.Item(455, 4).Value = "18.03.202570"
$ws.Cells.Item(456, 1).Value = 45734.72916666666
$ws.Cells.Item(456, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(456, 2).Value = 7780
$ws.Cells.Item(456, 3).Value = 71
$ws.Cells.Item(456, 4).Value = "18.03.202571"
$ws.Cells.Item(457, 1).Value = 45734.73958333334
$ws.Cells.Item(457, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(457, 2).Value = 7880
$ws.Cells.Item(457, 3).Value = 72
$ws.Cells.Item(457, 4).Value = "18.03.202572"
$ws.Cells.Item(458, 1).Value = 45734.75
$ws.Cells.Item(458, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(458, 2).Value = 7970
$ws.Cells.Item(458, 3).Value = 73
$ws.Cells.Item(458, 4).Value = "18.03.202573"
$ws.Cells.Item(459, 1).Value = 45734.76041666666
$ws.Cells.Item(459, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(459, 2).Value = 8040
$ws.Cells.Item(459, 3).Value = 74
$ws.Cells.Item(459, 4).Value = "18.03.202574"
$ws.Cells.Item(460, 1).Value = 45734.77083333334
$ws.Cells.Item(460, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(460, 2).Value = 8070
$ws.Cells.Item(460, 3).Value = 75
$ws.Cells.Item(460, 4).Value = "18.03.202575"
$ws.Cells.Item(461, 1).Value = 45734.78125
$ws.Cells.Item(461, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(461, 2).Value = 8060
$ws.Cells.Item(461, 3).Value = 76
$ws.Cells.Item(461, 4).Value = "18.03.202576"
$ws.Cells.Item(462, 1).Value = 45734.79166666666
$ws.Cells.Item(462, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(462, 2).Value = 8030
$ws.Cells.Item(462, 3).Value = 77
$ws.Cells.Item(462, 4).Value = "18.03.202577"
$ws.Cells.Item(463, 1).Value = 45734.80208333334
$ws.Cells.Item(463, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(463, 2).Value = 7980
$ws.Cells.Item(463, 3).Value = 78
$ws.Cells.Item(463, 4).Value = "18.03.202578"
$ws.Cells.Item(464, 1).Value = 45734.8125
$ws.Cells.Item(464, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(464, 2).Value = 7910
$ws.Cells.Item(464, 3).Value = 79
$ws.Cells.Item(464, 4).Value = "18.03.202579"
$ws.Cells.Item(465, 1).Value = 45734.82291666666
$ws.Cells.Item(465, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(465, 2).Value = 7800
$ws.Cells.Item(465, 3).Value = 80
$ws.Cells.Item(465, 4).Value = "18.03.202580"
$ws.Cells.Item(466, 1).Value = 45734.83333333334
$ws.Cells.Item(466, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(466, 2).Value = 7670
$ws.Cells.Item(466, 3).Value = 81
$ws.Cells.Item(466, 4).Value = "18.03.202581"
$ws.Cells.Item(467, 1).Value = 45734.84375
$ws.Cells.Item(467, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(467, 2).Value = 7550
$ws.Cells.Item(467, 3).Value = 82
$ws.Cells.Item(467, 4).Value = "18.03.202582"
$ws.Cells.Item(468, 1).Value = 45734.85416666666
$ws.Cells.Item(468, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(468, 2).Value = 7420
$ws.Cells.Item(468, 3).Value = 83
$ws.Cells.Item(468, 4).Value = "18.03.202583"
$ws.Cells.Item(469, 1).Value = 45734.86458333334
$ws.Cells.Item(469, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(469, 2).Value = 7300
$ws.Cells.Item(469, 3).Value = 84
$ws.Cells.Item(469, 4).Value = "18.03.202584"
$ws.Cells.Item(470, 1).Value = 45734.875
$ws.Cells.Item(470, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(470, 2).Value = 7170
$ws.Cells.Item(470, 3).Value = 85
$ws.Cells.Item(470, 4).Value = "18.03.202585"
$ws.Cells.Item(471, 1).Value = 45734.88541666666
$ws.Cells.Item(471, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(471, 2).Value = 7030
$ws.Cells.Item(471, 3).Value = 86
$ws.Cells.Item(471, 4).Value = "18.03.202586"
$ws.Cells.Item(472, 1).Value = 45734.89583333334
$ws.Cells.Item(472, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(472, 2).Value = 6900
$ws.Cells.Item(472, 3).Value = 87
$ws.Cells.Item(472, 4).Value = "18.03.202587"
$ws.Cells.Item(473, 1).Value = 45734.90625
$ws.Cells.Item(473, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(473, 2).Value = 6750
$ws.Cells.Item(473, 3).Value = 88
$ws.Cells.Item(473, 4).Value = "18.03.202588"
$ws.Cells.Item(474, 1).Value = 45734.91666666666
$ws.Cells.Item(474, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(474, 2).Value = 6620
$ws.Cells.Item(474, 3).Value = 89
$ws.Cells.Item(474, 4).Value = "18.03.202589"
$ws.Cells.Item(475, 1).Value = 45734.92708333334
$ws.Cells.Item(475, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(475, 2).Value = 6480
$ws.Cells.Item(475, 3).Value = 90
$ws.Cells.Item(475, 4).Value = "18.03.202590"
$ws.Cells.Item(476, 1).Value = 45734.9375
$ws.Cells.Item(476, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(476, 2).Value = 6340
$ws.Cells.Item(476, 3).Value = 91
$ws.Cells.Item(476, 4).Value = "18.03.202591"
$ws.Cells.Item(477, 1).Value = 45734.94791666666
$ws.Cells.Item(477, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(477, 2).Value = 6220
$ws.Cells.Item(477, 3).Value = 92
$ws.Cells.Item(477, 4).Value = "18.03.202592"
$ws.Cells.Item(478, 1).Value = 45734.95833333334
$ws.Cells.Item(478, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(478, 2).Value = 6020
$ws.Cells.Item(478, 3).Value = 93
$ws.Cells.Item(478, 4).Value = "18.03.202593"
$ws.Cells.Item(479, 1).Value = 45734.96875
$ws.Cells.Item(479, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(479, 2).Value = 5950
$ws.Cells.Item(479, 3).Value = 94
$ws.Cells.Item(479, 4).Value = "18.03.202594"
$ws.Cells.Item(480, 1).Value = 45734.97916666666
$ws.Cells.Item(480, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(480, 2).Value = 5920
$ws.Cells.Item(480, 3).Value = 95
$ws.Cells.Item(480, 4).Value = "18.03.202595"
$ws.Cells.Item(481, 1).Value = 45734.98958333334
$ws.Cells.Item(481, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(481, 2).Value = 5860
$ws.Cells.Item(481, 3).Value = 96
$ws.Cells.Item(481, 4).Value = "18.03.202596"

"Done: rows extended to 481"
